$d = $word.ActiveDocument

# 1) Merge the split " (Guller, 2015)" citation runs (removing proofErr spell-check wrapping)
$d.Content.Find.Execute(" (Guller, 2015)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " (Guller, 2015)", 2)

# 2) Add "Disadvantages" text to the empty paragraph following "Advantages"
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $p.Range.Text = "Disadvantages"
        break
    }
    if ($p.Range.Text.TrimEnd("`r`n").Trim() -eq "Advantages") {
        $found = $true
    }
}
